$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FinalDec")
$ws.Range("D:E").Delete()
